$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 66573
$ws.Range("E2").Value = 1311514368074
$ws.Range("F2").Value = 12152647814
$ws.Range("G2").Value = 0.48677
$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 3557.3
$ws.Range("E3").Value = 427234503583
$ws.Range("F3").Value = 9103359556
$ws.Range("G3").Value = 0.42272
$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.999292
$ws.Range("E4").Value = 112530262253
$ws.Range("F4").Value = 13895908664
$ws.Range("G4").Value = 0.00396
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 606.33
$ws.Range("E5").Value = 93238905955
$ws.Range("F5").Value = 577747781
$ws.Range("G5").Value = -0.376
$ws.Range("B6").Value = "SOL"
$ws.Range("C6").Value = "Solana"
$ws.Range("D6").Value = 145.31
$ws.Range("E6").Value = 67032872604
$ws.Range("F6").Value = 1101451423
$ws.Range("G6").Value = 0.72075
$ws.Range("B7").Value = "STETH"
$ws.Range("C7").Value = "Lido Staked Ether"
$ws.Range("D7").Value = 3555.75
$ws.Range("E7").Value = 33891275015
$ws.Range("F7").Value = 39252266
$ws.Range("G7").Value = 0.42802
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "USDC"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 32481212897
$ws.Range("F8").Value = 2347446916
$ws.Range("G8").Value = 0.07713
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "XRP"
$ws.Range("D9").Value = 0.495165
$ws.Range("E9").Value = 27528514562
$ws.Range("F9").Value = 757893422
$ws.Range("G9").Value = 3.28048
$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.135826
$ws.Range("E10").Value = 19644021810
$ws.Range("F10").Value = 385990823
$ws.Range("G10").Value = -0.56852
$ws.Range("B11").Value = "TON"
$ws.Range("C11").Value = "Toncoin"
$ws.Range("D11").Value = 7.99
$ws.Range("E11").Value = 19412338275
$ws.Range("F11").Value = 388797505
$ws.Range("G11").Value = -1.64512
$ws.Range("B12").Value = "ADA"
$ws.Range("C12").Value = "Cardano"
$ws.Range("D12").Value = 0.415052
$ws.Range("E12").Value = 14684055757
$ws.Range("F12").Value = 182088429
$ws.Range("G12").Value = 1.11794
$ws.Range("B13").Value = "SHIB"
$ws.Range("C13").Value = "Shiba Inu"
$ws.Range("D13").Value = 0.00002074
$ws.Range("E13").Value = 12207860611
$ws.Range("F13").Value = 172099522
$ws.Range("G13").Value = -0.06411
$ws.Range("B14").Value = "AVAX"
$ws.Range("C14").Value = "Avalanche"
$ws.Range("D14").Value = 30.03
$ws.Range("E14").Value = 11797111581
$ws.Range("F14").Value = 189082038
$ws.Range("G14").Value = -0.27576
$ws.Range("B15").Value = "WBTC"
$ws.Range("C15").Value = "Wrapped Bitcoin"
$ws.Range("D15").Value = 66597
$ws.Range("E15").Value = 10186043878
$ws.Range("F15").Value = 137996145
$ws.Range("G15").Value = 0.45414
$ws.Range("B16").Value = "TRX"
$ws.Range("C16").Value = "TRON"
$ws.Range("D16").Value = 0.115298
$ws.Range("E16").Value = 10058362159
$ws.Range("F16").Value = 196577834
$ws.Range("G16").Value = 0.10935
$ws.Range("B17").Value = "LINK"
$ws.Range("C17").Value = "Chainlink"
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 8789990886
$ws.Range("F17").Value = 237541729
$ws.Range("G17").Value = 0.56547
$ws.Range("B18").Value = "UNI"
$ws.Range("C18").Value = "Uniswap"
$ws.Range("D18").Value = 11.47
$ws.Range("E18").Value = 8642925175
$ws.Range("F18").Value = 381494663
$ws.Range("G18").Value = 4.79651
$ws.Range("B19").Value = "DOT"
$ws.Range("C19").Value = "Polkadot"
$ws.Range("D19").Value = 6.22
$ws.Range("E19").Value = 8554733323
$ws.Range("F19").Value = 98684591
$ws.Range("G19").Value = 0.26773
$ws.Range("B20").Value = "BCH"
$ws.Range("C20").Value = "Bitcoin Cash"
$ws.Range("D20").Value = 431.6
$ws.Range("E20").Value = 8506097314
$ws.Range("F20").Value = 127980838
$ws.Range("G20").Value = 1.30081
$ws.Range("B21").Value = "NEAR"
$ws.Range("C21").Value = "NEAR Protocol"
$ws.Range("D21").Value = 5.62
$ws.Range("E21").Value = 6116610822
$ws.Range("F21").Value = 153276288
$ws.Range("G21").Value = -0.19331
$ws.Range("B22").Value = "LTC"
$ws.Range("C22").Value = "Litecoin"
$ws.Range("D22").Value = 79.3
$ws.Range("E22").Value = 5916031717
$ws.Range("F22").Value = 217059682
$ws.Range("G22").Value = 0.40544
$ws.Range("B23").Value = "MATIC"
$ws.Range("C23").Value = "Polygon"
$ws.Range("D23").Value = 0.614164
$ws.Range("E23").Value = 5693888274
$ws.Range("F23").Value = 220684837
$ws.Range("G23").Value = 2.13205
$ws.Range("B24").Value = "WEETH"
$ws.Range("C24").Value = "Wrapped eETH"
$ws.Range("D24").Value = 3697.16
$ws.Range("E24").Value = 5442384333
$ws.Range("F24").Value = 19518612
$ws.Range("G24").Value = 0.53203
$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "LEO Token"
$ws.Range("D25").Value = 5.84
$ws.Range("E25").Value = 5408547350
$ws.Range("F25").Value = 976975
$ws.Range("G25").Value = 2.16434
$ws.Range("B26").Value = "DAI"
$ws.Range("C26").Value = "Dai"
$ws.Range("D26").Value = 0.998336
$ws.Range("E26").Value = 5234477104
$ws.Range("F26").Value = 309660990
$ws.Range("G26").Value = -0.06297999999999999
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "Pepe"
$ws.Range("D27").Value = 0.00001192
$ws.Range("E27").Value = 5002219868
$ws.Range("F27").Value = 498642594
$ws.Range("G27").Value = -0.92702
$ws.Range("B28").Value = "ICP"
$ws.Range("C28").Value = "Internet Computer"
$ws.Range("D28").Value = 9.16
$ws.Range("E28").Value = 4255774632
$ws.Range("F28").Value = 64719908
$ws.Range("G28").Value = -0.33152
$ws.Range("B29").Value = "ETC"
$ws.Range("C29").Value = "Ethereum Classic"
$ws.Range("D29").Value = 25.32
$ws.Range("E29").Value = 3732692357
$ws.Range("F29").Value = 77744818
$ws.Range("G29").Value = -0.04392
$ws.Range("B30").Value = "KAS"
$ws.Range("C30").Value = "Kaspa"
$ws.Range("D30").Value = 0.15431
$ws.Range("E30").Value = 3687473918
$ws.Range("F30").Value = 56796979
$ws.Range("G30").Value = -3.07067
$ws.Range("B31").Value = "FET"
$ws.Range("C31").Value = "Fetch.ai"
$ws.Range("D31").Value = 1.45
$ws.Range("E31").Value = 3665177792
$ws.Range("F31").Value = 87241720
$ws.Range("G31").Value = -1.90103
$ws.Range("B32").Value = "EZETH"
$ws.Range("C32").Value = "Renzo Restaked ETH"
$ws.Range("D32").Value = 3545.05
$ws.Range("E32").Value = 3532389876
$ws.Range("F32").Value = 38546355
$ws.Range("G32").Value = 0.63964
$ws.Range("B33").Value = "APT"
$ws.Range("C33").Value = "Aptos"
$ws.Range("D33").Value = 7.83
$ws.Range("E33").Value = 3527424973
$ws.Range("F33").Value = 125128556
$ws.Range("G33").Value = 0.23057
$ws.Range("B34").Value = "USDE"
$ws.Range("C34").Value = "Ethena USDe"
$ws.Range("D34").Value = 1.001
$ws.Range("E34").Value = 3524327678
$ws.Range("F34").Value = 28172087
$ws.Range("G34").Value = 0.08056000000000001
$ws.Range("B35").Value = "XMR"
$ws.Range("C35").Value = "Monero"
$ws.Range("D35").Value = 173.32
$ws.Range("E35").Value = 3195489421
$ws.Range("F35").Value = 54519453
$ws.Range("G35").Value = 0.20355
$ws.Range("B36").Value = "RNDR"
$ws.Range("C36").Value = "Render"
$ws.Range("D36").Value = 8.06
$ws.Range("E36").Value = 3128501298
$ws.Range("F36").Value = 73081803
$ws.Range("G36").Value = -0.51832
$ws.Range("B37").Value = "HBAR"
$ws.Range("C37").Value = "Hedera"
$ws.Range("D37").Value = 0.084943
$ws.Range("E37").Value = 3032834269
$ws.Range("F37").Value = 30866378
$ws.Range("G37").Value = -0.93716
$ws.Range("B38").Value = "FIL"
$ws.Range("C38").Value = "Filecoin"
$ws.Range("D38").Value = 5.19
$ws.Range("E38").Value = 2924493138
$ws.Range("F38").Value = 116811067
$ws.Range("G38").Value = -1.02652
$ws.Range("B39").Value = "MNT"
$ws.Range("C39").Value = "Mantle"
$ws.Range("D39").Value = 0.887687
$ws.Range("E39").Value = 2894677892
$ws.Range("F39").Value = 149708327
$ws.Range("G39").Value = -0.54756
$ws.Range("B40").Value = "XLM"
$ws.Range("C40").Value = "Stellar"
$ws.Range("D40").Value = 0.09858500000000001
$ws.Range("E40").Value = 2868242366
$ws.Range("F40").Value = 33744802
$ws.Range("G40").Value = 1.53861
$ws.Range("B41").Value = "STX"
$ws.Range("C41").Value = "Stacks"
$ws.Range("D41").Value = 1.94
$ws.Range("E41").Value = 2841483465
$ws.Range("F41").Value = 46703981
$ws.Range("G41").Value = 1.80814
$ws.Range("B42").Value = "ATOM"
$ws.Range("C42").Value = "Cosmos Hub"
$ws.Range("D42").Value = 7.17
$ws.Range("E42").Value = 2800475214
$ws.Range("F42").Value = 97099455
$ws.Range("G42").Value = 0.76178
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "OKB"
$ws.Range("D43").Value = 46.15
$ws.Range("E43").Value = 2764061492
$ws.Range("F43").Value = 2767354
$ws.Range("G43").Value = 1.54254
$ws.Range("B44").Value = "CRO"
$ws.Range("C44").Value = "Cronos"
$ws.Range("D44").Value = 0.099255
$ws.Range("E44").Value = 2657910274
$ws.Range("F44").Value = 7773194
$ws.Range("G44").Value = -2.63816
$ws.Range("B45").Value = "ARB"
$ws.Range("C45").Value = "Arbitrum"
$ws.Range("D45").Value = 0.917296
$ws.Range("E45").Value = 2656067608
$ws.Range("F45").Value = 120568164
$ws.Range("G45").Value = -0.57645
$ws.Range("B46").Value = "IMX"
$ws.Range("C46").Value = "Immutable"
$ws.Range("D46").Value = 1.72
$ws.Range("E46").Value = 2597179330
$ws.Range("F46").Value = 38097710
$ws.Range("G46").Value = -1.46766
$ws.Range("B47").Value = "FDUSD"
$ws.Range("C47").Value = "First Digital USD"
$ws.Range("D47").Value = 0.999582
$ws.Range("E47").Value = 2546129009
$ws.Range("F47").Value = 2663389695
$ws.Range("G47").Value = 0.03561
$ws.Range("B48").Value = "WIF"
$ws.Range("C48").Value = "dogwifhat"
$ws.Range("D48").Value = 2.52
$ws.Range("E48").Value = 2520038169
$ws.Range("F48").Value = 295054558
$ws.Range("G48").Value = 4.65933
$ws.Range("B49").Value = "INJ"
$ws.Range("C49").Value = "Injective"
$ws.Range("D49").Value = 25.09
$ws.Range("E49").Value = 2438375783
$ws.Range("F49").Value = 120958787
$ws.Range("G49").Value = -3.51505
$ws.Range("B50").Value = "SUI"
$ws.Range("C50").Value = "Sui"
$ws.Range("D50").Value = 0.944039
$ws.Range("E50").Value = 2288813310
$ws.Range("F50").Value = 99105152
$ws.Range("G50").Value = 0.15364
$ws.Range("B51").Value = "OP"
$ws.Range("C51").Value = "Optimism"
$ws.Range("D51").Value = 2.06
$ws.Range("E51").Value = 2238142443
$ws.Range("F51").Value = 154831016
$ws.Range("G51").Value = -0.11394
